$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Bmp7", "Bmpr2", "ECs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 36.81180933333333, 110.435428, 0.2598784967371026, 0.2598784967371026, 4.790333019149776, 43.112997172348, 0.07760301665407952, 0.07760301665407952),
    @("ECs", "Bmp7", "Bmpr2", "FAPs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 27.90901333333333, 83.72704, 0.1970278712683331, 0.197027871268333, 3.631809208071111, 32.68628287264, 0.05883502239441481, 0.0588350223944148),
    @("ECs", "Bmp7", "Bmpr2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 21.95609833333333, 65.868295, 0.1550023737603119, 0.1550023737603119, 2.857154394816111, 25.714389553345, 0.04628567558827974, 0.04628567558827974),
    @("ECs", "Bmp7", "Bmpr2", "MuSCs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 13.23098133333333, 39.692944, 0.09340609987756826, 0.09340609987756825, 1.721752011233777, 15.495768101104, 0.02789224662833241, 0.02789224662833241),
    @("ECs", "Bmp7", "Bmpr2", "Neutrophils", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 22.080681, 66.242043, 0.1558818838066577, 0.1558818838066577, 2.873366378756999, 25.860297408813, 0.04654830844798512, 0.04654830844798512),
    @("ECs", "Bmp7", "Bmpr2", "Resolving-Mac", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 19.66149466666667, 58.984484, 0.1388032745500265, 0.1388032745500265, 2.558556854804888, 23.027011693244, 0.04144841901807352, 0.04144841901807352),
    @("FAPs", "Bmp7", "Bmpr2", "ECs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 36.81180933333333, 110.435428, 0.2598784967371026, 0.2598784967371026, 10.46260336630755, 94.16343029676798, 0.1694933483820089, 0.1694933483820089),
    @("FAPs", "Bmp7", "Bmpr2", "FAPs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 27.90901333333333, 83.72704, 0.1970278712683331, 0.197027871268333, 7.932262557582223, 71.39036301824, 0.1285020270824178, 0.1285020270824177),
    @("FAPs", "Bmp7", "Bmpr2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 21.95609833333333, 65.868295, 0.1550023737603119, 0.1550023737603119, 6.240332993502222, 56.16299694152, 0.1010929017431249, 0.1010929017431248),
    @("FAPs", "Bmp7", "Bmpr2", "MuSCs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 13.23098133333333, 39.692944, 0.09340609987756826, 0.09340609987756825, 3.760491873251556, 33.844426859264, 0.06091967140924715, 0.06091967140924712),
    @("FAPs", "Bmp7", "Bmpr2", "Neutrophils", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 22.080681, 66.242043, 0.1558818838066577, 0.1558818838066577, 6.275741712912, 56.48167541620799, 0.1016665202016061, 0.1016665202016061),
    @("FAPs", "Bmp7", "Bmpr2", "Resolving-Mac", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 19.66149466666667, 58.984484, 0.1388032745500265, 0.1388032745500265, 5.588163798833778, 50.293474189504, 0.09052781228633472, 0.09052781228633469),
    @("MuSCs", "Bmp7", "Bmpr2", "ECs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 36.81180933333333, 110.435428, 0.2598784967371026, 0.2598784967371026, 0.7890243212506665, 7.101218891255999, 0.01278213170101417, 0.01278213170101417),
    @("MuSCs", "Bmp7", "Bmpr2", "FAPs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 27.90901333333333, 83.72704, 0.1970278712683331, 0.197027871268333, 0.5982017917866667, 5.38381612608, 0.009690821791500473, 0.009690821791500472),
    @("MuSCs", "Bmp7", "Bmpr2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 21.95609833333333, 65.868295, 0.1550023737603119, 0.1550023737603119, 0.4706070116766666, 4.23546310509, 0.007623796428907336, 0.007623796428907336),
    @("MuSCs", "Bmp7", "Bmpr2", "MuSCs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 13.23098133333333, 39.692944, 0.09340609987756826, 0.09340609987756825, 0.2835928538986666, 2.552335685088, 0.004594181839988705, 0.004594181839988705),
    @("MuSCs", "Bmp7", "Bmpr2", "Neutrophils", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 22.080681, 66.242043, 0.1558818838066577, 0.1558818838066577, 0.4732773165539999, 4.259495848986, 0.007667055157066479, 0.007667055157066479),
    @("MuSCs", "Bmp7", "Bmpr2", "Resolving-Mac", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 19.66149466666667, 58.984484, 0.1388032745500265, 0.1388032745500265, 0.4214244766853333, 3.792820290168, 0.006827043245618273, 0.006827043245618273)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
